$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DEV 7 - Manager Bug Fixes: project "Acacia Breeze" (row 2) should be
# Hidden, not Visible.
$ws.Range("O2").Value = "Hidden"
